$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("L1").Value = "MPIThreadSafe"
$ws.Range("M1").Value = "MPISendRec"
$ws.Range("J1").Value = "MPISeq"
$ws.Range("K1").Value = "MPIParBaseNodes"
$ws.Range("G1").Value = "ParBase"
$ws.Range("H1").Value = "MPIParBase"

# --- Row 2 ---
$ws.Range("G2").Value = 331.12
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 327.678245
$ws.Range("K2").Value = 330.709878
$ws.Range("L2").Value = 330.709878
$ws.Range("M2").Value = 327.678245

# --- Row 3 ---
$ws.Range("G3").Value = 331.269969
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 175.07313300000001
$ws.Range("K3").Value = 172.31019599999999
$ws.Range("L3").Value = 172.460532
$ws.Range("M3").Value = 177.07665700000001

# --- Row 4 ---
$ws.Range("H4").Value = 90.144876999999994
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 116.310136
$ws.Range("K4").Value = 117.08775799999999
$ws.Range("L4").Value = 117.18559999999999
$ws.Range("M4").Value = 119.907118

# --- Row 5 ---
$ws.Range("H5").Value = 90.355018000000001
$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 88.853880000000004
$ws.Range("K5").Value = 89.760159999999999
$ws.Range("L5").Value = 89.774897999999993
$ws.Range("M5").Value = 91.457217999999997

# --- Column widths for newly used / resized columns ---
$ws.Columns.Item(11).ColumnWidth = 17.42578125
$ws.Columns.Item(12).ColumnWidth = 14.5703125
$ws.Columns.Item(13).ColumnWidth = 12

# --- Selection state to match saved view ---
$ws.Range("K5").Select()
